$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 10; $r++) {
    $ws.Range("C$r").Value = 45184
}
